$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 75 ---------------------------------------------------------------
# Column A holds a date-like string ("2024-05-15"). Assigning it straight to
# .Value would make Excel auto-convert it to a real date serial number, which
# is not what the source file has (it is stored as plain text). Writing it as
# a formula that evaluates to the literal text, then copy/paste-special as
# values, bakes in a static text cell without Excel's autoformat kicking in.
$ws.Range("A75").Formula = "=""2024-05-15"""
$ws.Range("A75").Copy()
$ws.Range("A75").PasteSpecial(-4163)

$ws.Range("B75").Value = "13:00:36"
$ws.Range("C75").Value = "Palet atascado en la curva"
$ws.Range("D75").Value = "-"
$ws.Range("E75").Value = "-"
$ws.Range("F75").Value = "-"
$ws.Range("G75").Value = "-"

# --- Row 76 ---------------------------------------------------------------
$ws.Range("A76").Formula = "=""2024-05-15"""
$ws.Range("A76").Copy()
$ws.Range("A76").PasteSpecial(-4163)

$ws.Range("B76").Value = "13:00:42"
$ws.Range("C76").Value = "Fallo en elevador"
$ws.Range("D76").Value = "-"
$ws.Range("E76").Value = "-"
$ws.Range("F76").Value = "-"
$ws.Range("G76").Value = "-"

# H76 is an empty but present text cell in the source. A quote-prefix forces
# text type on an otherwise-empty value; resetting the style back to Normal
# afterwards drops the quote-prefix formatting again while keeping the cell
# typed as text.
$ws.Range("H76").Value = "'"
$ws.Range("H76").Style = "Normal"

$excel.CutCopyMode = $false
